$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text cells keep their literal string representation (avoid Excel
# auto-converting number-like or percent-like text into real numbers/dates).
$textCells = @('D2', 'E2', 'D3', 'E3', 'E4', 'D5', 'E5', 'D6', 'E6', 'D7', 'E7', 'E8', 'E9', 'E10', 'D11', 'E11', 'D12', 'E12', 'D13', 'E13', 'D14', 'E14', 'D15', 'E15', 'D16', 'E16', 'D17', 'E17', 'E18', 'D19', 'E19', 'D20', 'E20', 'D21', 'E21', 'D22', 'E22', 'D23', 'E23', 'D24', 'E24', 'D25', 'E25', 'E26', 'E27', 'D28', 'E28', 'D29', 'E29', 'E30', 'E31', 'E32', 'D33', 'E33', 'D34', 'E34', 'D35', 'E35', 'E36', 'D37', 'E37', 'D38', 'E38', 'D39', 'E39', 'E40', 'D41', 'E41', 'D42', 'E42', 'D43', 'E43', 'D44', 'E44', 'E45', 'D46', 'E46', 'B47', 'C47', 'D47', 'E47', 'B48', 'C48', 'D48', 'E48', 'E49', 'E50', 'D51', 'E51')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '66.183.04'
$ws.Range('E2').Value = '  -1.14%  '
$ws.Range('D3').Value = '3.531.64'
$ws.Range('E3').Value = '  +0.46%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').Value = '607.16'
$ws.Range('E5').Value = '  -0.17%  '
$ws.Range('D6').Value = '143.50'
$ws.Range('E6').Value = '  -2.94%  '
$ws.Range('D7').Value = '3.531.89'
$ws.Range('E7').Value = '  +0.48%  '
$ws.Range('E8').Value = '  -0.08%  '
$ws.Range('E9').Value = '  +0.32%  '
$ws.Range('E10').Value = '  -4.26%  '
$ws.Range('D11').Value = '8.04'
$ws.Range('E11').Value = '  +1.26%  '
$ws.Range('D12').Value = '0.411'
$ws.Range('E12').Value = '  -2.70%  '
$ws.Range('D13').Value = '4.126.00'
$ws.Range('E13').Value = '  +0.39%  '
$ws.Range('D14').Value = '0.0000207'
$ws.Range('E14').Value = '  -4.71%  '
$ws.Range('D15').Value = '30.27'
$ws.Range('E15').Value = '  -5.19%  '
$ws.Range('D16').Value = '3.531.12'
$ws.Range('E16').Value = '  +0.25%  '
$ws.Range('D17').Value = '66.266.63'
$ws.Range('E17').Value = '  -1.07%  '
$ws.Range('E18').Value = '  -0.60%  '
$ws.Range('D19').Value = '10.88'
$ws.Range('E19').Value = '  +1.70%  '
$ws.Range('D20').Value = '6.21'
$ws.Range('E20').Value = '  -3.36%  '
$ws.Range('D21').Value = '14.94'
$ws.Range('E21').Value = '  -2.47%  '
$ws.Range('D22').Value = '425.95'
$ws.Range('E22').Value = '  -2.74%  '
$ws.Range('D23').Value = '0.601'
$ws.Range('E23').Value = '  -1.30%  '
$ws.Range('D24').Value = '78.69'
$ws.Range('E24').Value = '  -0.67%  '
$ws.Range('D25').Value = '3.667.40'
$ws.Range('E25').Value = '  +0.25%  '
$ws.Range('E26').Value = '  -0.01%  '
$ws.Range('E27').Value = '  -1.93%  '
$ws.Range('D28').Value = '8.01'
$ws.Range('E28').Value = '  -3.43%  '
$ws.Range('D29').Value = '9.19'
$ws.Range('E29').Value = '  -5.88%  '
$ws.Range('E30').Value = '  -1.78%  '
$ws.Range('E31').Value = '  +0.10%  '
$ws.Range('E32').Value = '  -3.14%  '
$ws.Range('D33').Value = '1.48'
$ws.Range('E33').Value = '  -6.19%  '
$ws.Range('D34').Value = '25.29'
$ws.Range('E34').Value = '  -0.69%  '
$ws.Range('D35').Value = '3.518.59'
$ws.Range('E35').Value = '  +0.25%  '
$ws.Range('E36').Value = '  -0.03%  '
$ws.Range('D37').Value = '1.75'
$ws.Range('E37').Value = '  -2.82%  '
$ws.Range('D38').Value = '7.84'
$ws.Range('E38').Value = '  -2.19%  '
$ws.Range('D39').Value = '5.61'
$ws.Range('E39').Value = '  -5.87%  '
$ws.Range('E40').Value = '  -0.07%  '
$ws.Range('D41').Value = '172.29'
$ws.Range('E41').Value = '  -0.45%  '
$ws.Range('D42').Value = '0.0856'
$ws.Range('E42').Value = '  -4.20%  '
$ws.Range('D43').Value = '5.18'
$ws.Range('E43').Value = '  -4.39%  '
$ws.Range('D44').Value = '0.893'
$ws.Range('E44').Value = '  -0.38%  '
$ws.Range('E45').Value = '  -9.18%  '
$ws.Range('D46').Value = '45.34'
$ws.Range('E46').Value = '  -1.52%  '
$ws.Range('B47').Value = 'ONDO'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D47').Value = '1.21'
$ws.Range('E47').Value = '  -6.30%  '
$ws.Range('B48').Value = 'InjectiveProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D48').Value = '25.88'
$ws.Range('E48').Value = '  -7.90%  '
$ws.Range('E49').Value = '  -2.12%  '
$ws.Range('E50').Value = '  -4.05%  '
$ws.Range('D51').Value = '0.945'
$ws.Range('E51').Value = '  -4.30%  '
